$d = $word.ActiveDocument

# Locate the existing list paragraph that should follow our new bullet:
# "Please examine one or more scenarios ..."
$found = $d.Content
$found.Find.Execute("What do you think of the hypermedia aspects of this API and the use of HAL?", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$prevPara = $found.Paragraphs(1)

# Collapse to a point just before that paragraph's ending mark, so the new
# paragraph we insert lands between it and the next ("Please examine ...") paragraph.
$insertionPoint = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End - 1)

# Build the new list-item paragraph as raw OOXML (same numbering as its siblings:
# numId 2 / ilvl 0) with a trailing manual line break run, matching the
# formatting already used by every other bullet in this list.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:pPr><w:numPr><w:numId w:val="2"/><w:ilvl w:val="0"/></w:numPr></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">What do you think of our pagination scheme?</w:t></w:r>' + `
    '<w:r><w:cr/></w:r>' + `
    '</w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($packageXml)
